# Fala-Cidade_Apresentacao.pptx — "feat: correcao da apresentacao"
#
# 1) Remove the duplicate "Arquitetura do Sistema: Diagrama de Classes"
#    slide (slide 10 of 12) — it was an accidental duplicate of slide 9.
# 2) Tidy up the wording on the "Objetivo do Projeto" (Objetivos Não
#    Funcionais) slide: the bullet about the public-management tool had
#    been split across three runs ("...para " / "receber, e " /
#    "gerenciar as denúncias."); merge it back into a single run with
#    the same final text.

$p = $ppt.ActivePresentation

# --- 1) Delete the duplicate "Diagrama de Classes" slide ---------------
$p.Slides.Item(10).Delete()

# --- 2) Clean up the split run on the objectives slide ------------------
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange
$target = [char]0x0046 + "ornecer " + [char]0x00E0 + " gest" + [char]0x00E3 + `
    "o p" + [char]0x00FA + "blica uma ferramenta para receber, e gerenciar as den" + `
    [char]0x00FA + "ncias."
$para = $tr.Paragraphs(5, 1)

# Force a real text mutation so the engine rebuilds the paragraph as a
# single run (an assignment that is character-for-character identical to
# the existing text is treated as a no-op and would leave the original
# run split untouched).
$para.Text = "placeholder"
$para.Text = $target
